$p = $ppt.ActivePresentation

# The existing last slide (slide 5) is "The Casual heroes..." -- the team
# bios slide. The edit: duplicate it to become a new final slide (with the
# team list reordered), and repurpose the original slide 5 into a new
# "The Technology..." slide describing the technology stack.

$heroesSlide = $p.Slides.Item(5)

# 1) Duplicate slide 5 -> new slide 6, keeps "The Casual heroes..." title
#    and picture, but the team members get reordered on the new slide.
$dupRange = $heroesSlide.Duplicate()
$newSlide = $dupRange.Item(1)

$newBody = $newSlide.Shapes.Item(2).TextFrame.TextRange
$newBody.Text = "William Wedin`rTor Harrington`rRobert McCartney`rPatrick Mooney`rJoel Hammond-Turner`rJames Counihan"

# Re-split "William " / "Wedin" and "James " / "Counihan" into separate
# runs (mirrors the misspelling-flagged runs from the source deck).
$p1 = $newBody.Characters(1, 8)
$p1.Text = "William "
$p2 = $newBody.Characters(9, 5)
$p2.Text = "Wedin"

$fullLen = $newBody.Length
$jamesStart = $fullLen - [string]"Counihan".Length - [string]"James ".Length + 1
$jPart = $newBody.Characters($jamesStart, 6)
$jPart.Text = "James "
$cPart = $newBody.Characters($jamesStart + 6, 8)
$cPart.Text = "Counihan"

# 2) Turn the original slide 5 into "The Technology..." slide.
$title = $heroesSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "The Technology…"

$body = $heroesSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "HTML5 / jQuery`rREST WebAPI Services`rAzure Mobile Services`rSQL Azure Database`rGithub source control`rAzure Websites`rGoogle Maps`rNative Android client / Java"

# Split "HTML5 / " / "jQuery" into separate runs.
$r1 = $body.Characters(1, 8)
$r1.Text = "HTML5 / "
$r2 = $body.Characters(9, 6)
$r2.Text = "jQuery"

# Split "REST " / "WebAPI" / " Services" into separate runs.
$parasText = $body.Text
# Recompute offsets dynamically rather than hard-coding, to stay correct
# regardless of prior edits.
